$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2, E2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.923.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("E2").Style = "Normal"

# Row 3: D3, E3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.669.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E3").Style = "Normal"

# Row 4: E4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E4").Style = "Normal"

# Row 5: D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("E5").Style = "Normal"

# Row 6: D6, E6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.517"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("E6").Style = "Normal"

# Row 7: E7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E7").Style = "Normal"

# Row 8: E8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("E8").Style = "Normal"

# Row 9: E9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("E9").Style = "Normal"

# Row 10: D10, E10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("E10").Style = "Normal"

# Row 11: D11, E11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0893"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.75%  "
$ws.Range("E11").Style = "Normal"

# Row 12: D12, E12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.904.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("E12").Style = "Normal"

# Row 13: D13, E13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.667.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("E13").Style = "Normal"

# Row 14: E14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("E14").Style = "Normal"

# Row 15: D15, E15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.528"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("E15").Style = "Normal"

# Row 16: D16, E16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("E16").Style = "Normal"

# Row 17: D17, E17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.910.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("E17").Style = "Normal"

# Row 18: D18, E18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "235.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("E18").Style = "Normal"

# Row 19: E19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("E19").Style = "Normal"

# Row 20: D20, E20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0735"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("E20").Style = "Normal"

# Row 21: E21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E21").Style = "Normal"

# Row 22: D22, E22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("E22").Style = "Normal"

# Row 23: E23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("E23").Style = "Normal"

# Row 24: D24, E24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.51%  "
$ws.Range("E24").Style = "Normal"

# Row 25: D25, E25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("E25").Style = "Normal"

# Row 26: D26, E26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E26").Style = "Normal"

# Row 27: B27, C27, D27, E27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E27").Style = "Normal"

# Row 28: B28, C28, D28, E28
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.112"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("E28").Style = "Normal"

# Row 29: E29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E29").Style = "Normal"

# Row 30: D30, E30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0497"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("E30").Style = "Normal"

# Row 31: E31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E31").Style = "Normal"

# Row 32: E32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("E32").Style = "Normal"

# Row 33: D33, E33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.454.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.10%  "
$ws.Range("E33").Style = "Normal"

# Row 34: E34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.79%  "
$ws.Range("E34").Style = "Normal"

# Row 35: E35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.69%  "
$ws.Range("E35").Style = "Normal"

# Row 36: E36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E36").Style = "Normal"

# Row 37: D37, E37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.590"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.94%  "
$ws.Range("E37").Style = "Normal"

# Row 38: D38, E38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.906"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("E38").Style = "Normal"

# Row 39: E39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("E39").Style = "Normal"

# Row 40: D40, E40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.66%  "
$ws.Range("E40").Style = "Normal"

# Row 41: E41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E41").Style = "Normal"

# Row 42: E42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("E42").Style = "Normal"

# Row 43: D43, E43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.992"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +8.12%  "
$ws.Range("E43").Style = "Normal"

# Row 44: D44, E44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("E44").Style = "Normal"

# Row 45: D45, E45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.810.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("E45").Style = "Normal"

# Row 46: D46, E46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.783"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("E46").Style = "Normal"

# Row 47: D47, E47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("E47").Style = "Normal"

# Row 48: E48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("E48").Style = "Normal"

# Row 49: D49, E49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.102"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.32%  "
$ws.Range("E49").Style = "Normal"

# Row 50: E50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E50").Style = "Normal"

# Row 51: D51, E51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.34%  "
$ws.Range("E51").Style = "Normal"
